# Insert a new column before column A to hold row identifiers ("ID").
# This shifts existing columns A:E -> B:F.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").EntireColumn.Insert()

# Header for the new ID column - match style/formatting of the other header cells.
$ws.Range("B1").Copy($ws.Range("A1"))
$ws.Range("A1").Value = "ID"

# Row identifiers for rows 2-25 (data rows), now living in column A.
$ids = @(
    "Hb 2", "Hb 3", "S 24", "S 28", "Hb 107", "Hb 66", "Hb 69", "Hb 95",
    "Hb 99", "Hb 92", "Hb 40", "Hb 41", "S 11", "Hb 57", "S 21", "S 22",
    "S 3", "S 4", "S 5", "Hb 74", "Hb 79", "Hb 32", "S 15", "S 16"
)

$row = 2
foreach ($id in $ids) {
    $ws.Cells.Item($row, 1).Value = $id
    $row = $row + 1
}
